$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows (2 and 5) for the new "Alterno" scenario values ---

# Row 2: usuario autotest25 -> autotest30, valorRecarga 1000000 -> 160000
$ws.Cells.Item(2, 4).Value = "autotest30"
$ws.Cells.Item(2, 13).Value = "160000"

# Row 5: valorRecarga 25000 -> 30000, numeroCuenta 406-733020-17 -> 406-733020-18
$ws.Cells.Item(5, 13).Value = "30000"
$ws.Cells.Item(5, 15).Value = "406-733020-18"

# --- Add new row 6, cloning row 5's formatting first ---
$ws.Range("A5:O5").Copy()
$ws.Range("A6:O6").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(6).RowHeight = $ws.Rows.Item(5).RowHeight

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 93221453
$ws.Cells.Item(6, 3).Value = 1
$ws.Cells.Item(6, 4).Value = "autotest30"
$ws.Cells.Item(6, 5).Value = 1234
$ws.Cells.Item(6, 6).Value = 4321
$ws.Cells.Item(6, 7).Value = "Alterno"
$ws.Cells.Item(6, 8).NumberFormat = "@"
$ws.Cells.Item(6, 8).Value = "000"
$ws.Cells.Item(6, 9).NumberFormat = "@"
$ws.Cells.Item(6, 9).Value = "0369"
$ws.Cells.Item(6, 10).Value = "NO ERROR"
$ws.Cells.Item(6, 11).Value = "bolp"
$ws.Cells.Item(6, 12).Value = "ACTIVO"
$ws.Cells.Item(6, 13).NumberFormat = "@"
$ws.Cells.Item(6, 13).Value = "100000"
$ws.Cells.Item(6, 14).Value = "Ahorros"
$ws.Cells.Item(6, 15).Value = "406-733020-19"

$ws.Range("M6").Select()
